$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("L1").Value = "Excused Folio Ids"
